# Eurostat AT cars road eqr carpda: roll the TIME_PERIOD window forward by one
# year (drop 2012, each subsequent year's data shifts left, 2023 becomes the
# new rightmost column) and fix a label typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dropping the whole F column (year 2012, which was entirely blank) shifts
# every later column one position to the left - G2018->F2018... 2013 header
# moves to F1, ..., 2023 header moves to P1 - exactly matching "the oldest
# year column is retired, the data set gains 2023 at the end".
$ws.Columns("F").Delete()

# Fix capitalisation typo in the Alternative Energy label.
$ws.Range("D7").Value = "Alternative energy [ALT]"
